$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.316.81"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "3.326.92"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").Value = "3.908.58"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "66.457.69"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "3.312.23"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "433.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.516"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000115"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.193"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.20%  "
$ws.Range("D39").Value = "2.813.42"
$ws.Range("E39").Value = "  +4.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "324.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0271"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.969"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.51%  "
